$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "have to login the first time he/she uses the application" is
# split in two runs with a collapsed "_GoBack" bookmark in between (moving
# the bookmark that used to sit near "request is submitted").
# ---------------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute("have to l", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r1.Collapse(0)
$d.Bookmarks.Add("_GoBack", $r1)

# ---------------------------------------------------------------------------
# Change 2: "Then the " -> "Next time" + " the " (keep neighbouring runs
# untouched by isolating the word "Then" with temporary bookmarks before
# replacing its text, then removing the temporary bookmarks again).
# ---------------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute("Then", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$thenStart = $r2.Start
$thenEnd = $r2.End

$tmpStart = $d.Range($thenStart, $thenStart)
$d.Bookmarks.Add("ZZ_TMP_START", $tmpStart)
$tmpEnd = $d.Range($thenEnd, $thenEnd)
$d.Bookmarks.Add("ZZ_TMP_END", $tmpEnd)

$thenRange = $d.Range($thenStart, $thenEnd)
$thenRange.Text = "Next time"

$d.Bookmarks("ZZ_TMP_START").Delete()
$d.Bookmarks("ZZ_TMP_END").Delete()

# ---------------------------------------------------------------------------
# Change 3: Replace the "Register" paragraph description text entirely.
# ---------------------------------------------------------------------------
$r3 = $d.Content
$r3.Find.Execute("If the user does not have credentials, he/she must register through the register user form. ", $false, $false, $false, $false, $false, $true, 1, $false, "Creates new user.", 2)

# ---------------------------------------------------------------------------
# Change 4: "/her forgotten " -> "/her " (keep neighbouring runs untouched).
# ---------------------------------------------------------------------------
$r4 = $d.Content
$r4.Find.Execute("/her forgotten ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r4Start = $r4.Start
$r4End = $r4.End

$tmp4Start = $d.Range($r4Start, $r4Start)
$d.Bookmarks.Add("ZZ_TMP_START2", $tmp4Start)
$tmp4End = $d.Range($r4End, $r4End)
$d.Bookmarks.Add("ZZ_TMP_END2", $tmp4End)

$forgottenRange = $d.Range($r4Start, $r4End)
$forgottenRange.Text = "/her "

$d.Bookmarks("ZZ_TMP_START2").Delete()
$d.Bookmarks("ZZ_TMP_END2").Delete()

# ---------------------------------------------------------------------------
# Change 5: Mark every inline picture's run as NoProof (adds
# <w:rPr><w:noProof/></w:rPr> ahead of each <w:drawing>).
# ---------------------------------------------------------------------------
$shapes = $d.InlineShapes
for ($i = 1; $i -le $shapes.Count; $i++) {
    $shapeRange = $shapes.Item($i).Range
    $shapeRange.NoProofing = 1
}

# ---------------------------------------------------------------------------
# Change 6: Merge the two runs around the (now relocated) old "_GoBack"
# bookmark location back into a single run: "...where the request is
# submitted."
# ---------------------------------------------------------------------------
$r6 = $d.Content
$r6.Find.Execute("where the req", $false, $false, $false, $false, $false, $true, 1, $false, "where the req", 2)

Write-Host "Done"
